# Auto-generated edit script: update Sheets via scheduled runner
# Applies numeric updates to market-price / profit columns (H-N) across
# several leve rows in the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 18
$ws_ALC.Range("H18").Value = 340.4
$ws_ALC.Range("I18").Value = 300
$ws_ALC.Range("K18").Value = 300
$ws_ALC.Range("M18").Value = -16

# ALC row 33
$ws_ALC.Range("H33").Value = 1036.4445
$ws_ALC.Range("I33").Value = 1022.7059
$ws_ALC.Range("J33").Value = 1059.8
$ws_ALC.Range("K33").Value = 1022.7059
$ws_ALC.Range("L33").Value = 1059.8
$ws_ALC.Range("M33").Value = -793.7059
$ws_ALC.Range("N33").Value = -1517.8

# ALC row 69
$ws_ALC.Range("H69").Value = 2993.1667
$ws_ALC.Range("I69").Value = 2653.3333
$ws_ALC.Range("J69").Value = 3333
$ws_ALC.Range("K69").Value = 7959.999899999999
$ws_ALC.Range("L69").Value = 9999
$ws_ALC.Range("M69").Value = -7085.999899999999
$ws_ALC.Range("N69").Value = -11747

# ALC row 72
$ws_ALC.Range("H72").Value = 2993.1667
$ws_ALC.Range("I72").Value = 2653.3333
$ws_ALC.Range("J72").Value = 3333
$ws_ALC.Range("K72").Value = 23879.9997
$ws_ALC.Range("L72").Value = 29997
$ws_ALC.Range("M72").Value = -19511.9997
$ws_ALC.Range("N72").Value = -38733

# ALC row 114
$ws_ALC.Range("H114").Value = 52180.5
$ws_ALC.Range("J114").Value = 52180.5
$ws_ALC.Range("L114").Value = 52180.5
$ws_ALC.Range("N114").Value = -60858.5

# ALC row 132
$ws_ALC.Range("H132").Value = 2515.3684
$ws_ALC.Range("I132").Value = 2373.5833
$ws_ALC.Range("J132").Value = 5067.5
$ws_ALC.Range("K132").Value = 7120.749899999999
$ws_ALC.Range("L132").Value = 15202.5
$ws_ALC.Range("M132").Value = -4590.749899999999
$ws_ALC.Range("N132").Value = -20262.5

# ALC row 137
$ws_ALC.Range("H137").Value = 3687.0688
$ws_ALC.Range("I137").Value = 3443.4546
$ws_ALC.Range("J137").Value = 4452.7144
$ws_ALC.Range("K137").Value = 10330.3638
$ws_ALC.Range("L137").Value = 13358.1432
$ws_ALC.Range("M137").Value = -7780.363799999999
$ws_ALC.Range("N137").Value = -18458.1432

# ALC row 138
$ws_ALC.Range("H138").Value = 186232.92
$ws_ALC.Range("I138").Value = 2169.524
$ws_ALC.Range("K138").Value = 6508.572
$ws_ALC.Range("M138").Value = -1368.572

# ARM row 32
$ws_ARM.Range("H32").Value = 409819.53
$ws_ARM.Range("I32").Value = 523033.53
$ws_ARM.Range("J32").Value = 10240.647
$ws_ARM.Range("K32").Value = 523033.53
$ws_ARM.Range("L32").Value = 10240.647
$ws_ARM.Range("M32").Value = -522746.53
$ws_ARM.Range("N32").Value = -10814.647

# ARM row 61
$ws_ARM.Range("H61").Value = 2495.25
$ws_ARM.Range("I61").Value = 1993.6666
$ws_ARM.Range("J61").Value = 4000
$ws_ARM.Range("K61").Value = 1993.6666
$ws_ARM.Range("L61").Value = 4000
$ws_ARM.Range("M61").Value = -1781.6666
$ws_ARM.Range("N61").Value = -4424

# ARM row 132
$ws_ARM.Range("H132").Value = 4818.65
$ws_ARM.Range("I132").Value = 5672.4
$ws_ARM.Range("J132").Value = 3964.9
$ws_ARM.Range("K132").Value = 17017.2
$ws_ARM.Range("L132").Value = 11894.7
$ws_ARM.Range("M132").Value = -14487.2
$ws_ARM.Range("N132").Value = -16954.7

# ARM row 136
$ws_ARM.Range("H136").Value = 2495.25
$ws_ARM.Range("I136").Value = 1993.6666
$ws_ARM.Range("J136").Value = 4000
$ws_ARM.Range("K136").Value = 5980.9998
$ws_ARM.Range("L136").Value = 12000
$ws_ARM.Range("M136").Value = -3430.9998
$ws_ARM.Range("N136").Value = -17100

# CRP row 58
$ws_CRP.Range("H58").Value = 1193.4166
$ws_CRP.Range("I58").Value = 829.8333
$ws_CRP.Range("J58").Value = 1557
$ws_CRP.Range("K58").Value = 829.8333
$ws_CRP.Range("L58").Value = 1557
$ws_CRP.Range("M58").Value = -626.8333
$ws_CRP.Range("N58").Value = -1963

# CRP row 94
$ws_CRP.Range("H94").Value = 2190
$ws_CRP.Range("J94").Value = 2190
$ws_CRP.Range("L94").Value = 2190
$ws_CRP.Range("N94").Value = -3092

# CRP row 134
$ws_CRP.Range("H134").Value = 1641.75
$ws_CRP.Range("I134").Value = 1447.7142
$ws_CRP.Range("K134").Value = 4343.142599999999
$ws_CRP.Range("M134").Value = -1808.142599999999

# CRP row 136
$ws_CRP.Range("H136").Value = 1193.4166
$ws_CRP.Range("I136").Value = 829.8333
$ws_CRP.Range("J136").Value = 1557
$ws_CRP.Range("K136").Value = 2489.4999
$ws_CRP.Range("L136").Value = 4671
$ws_CRP.Range("M136").Value = 60.5001000000002
$ws_CRP.Range("N136").Value = -9771

# CUL row 4
$ws_CUL.Range("H4").Value = 10003329
$ws_CUL.Range("I4").Value = 8004600
$ws_CUL.Range("J4").Value = 11113733
$ws_CUL.Range("K4").Value = 24013800
$ws_CUL.Range("L4").Value = 33341199
$ws_CUL.Range("M4").Value = -24013688
$ws_CUL.Range("N4").Value = -33341423

# CUL row 6
$ws_CUL.Range("H6").Value = 148.71428
$ws_CUL.Range("I6").Value = 56.833332
$ws_CUL.Range("J6").Value = 700
$ws_CUL.Range("K6").Value = 170.499996
$ws_CUL.Range("L6").Value = 2100
$ws_CUL.Range("M6").Value = -57.49999600000001
$ws_CUL.Range("N6").Value = -2326

# GSM row 122
$ws_GSM.Range("H122").Value = 5822.353
$ws_GSM.Range("I122").Value = 5000
$ws_GSM.Range("J122").Value = 5873.75
$ws_GSM.Range("K122").Value = 15000
$ws_GSM.Range("L122").Value = 17621.25
$ws_GSM.Range("M122").Value = -12550
$ws_GSM.Range("N122").Value = -22521.25

# GSM row 126
$ws_GSM.Range("H126").Value = 2584.9092
$ws_GSM.Range("J126").Value = 2417.4285
$ws_GSM.Range("L126").Value = 7252.2855
$ws_GSM.Range("N126").Value = -12192.2855

# GSM row 132
$ws_GSM.Range("H132").Value = 3608.5
$ws_GSM.Range("I132").Value = 2000
$ws_GSM.Range("J132").Value = 3838.2856
$ws_GSM.Range("K132").Value = 6000
$ws_GSM.Range("L132").Value = 11514.8568
$ws_GSM.Range("M132").Value = -3470
$ws_GSM.Range("N132").Value = -16574.8568

# LTW row 68
$ws_LTW.Range("H68").Value = 3537.634
$ws_LTW.Range("I68").Value = 2324.95
$ws_LTW.Range("J68").Value = 4692.5713
$ws_LTW.Range("K68").Value = 2324.95
$ws_LTW.Range("L68").Value = 4692.5713
$ws_LTW.Range("M68").Value = -1575.95
$ws_LTW.Range("N68").Value = -6190.5713

# LTW row 71
$ws_LTW.Range("H71").Value = 3537.634
$ws_LTW.Range("I71").Value = 2324.95
$ws_LTW.Range("J71").Value = 4692.5713
$ws_LTW.Range("K71").Value = 11624.75
$ws_LTW.Range("L71").Value = 23462.8565
$ws_LTW.Range("M71").Value = -7880.75
$ws_LTW.Range("N71").Value = -30950.8565

# LTW row 122
$ws_LTW.Range("H122").Value = 3000
$ws_LTW.Range("I122").Value = 0
$ws_LTW.Range("J122").Value = 3000
$ws_LTW.Range("K122").Value = 0
$ws_LTW.Range("L122").Value = 9000
$ws_LTW.Range("M122").ClearContents()
$ws_LTW.Range("N122").Value = -13900

# LTW row 132
$ws_LTW.Range("H132").Value = 3756.3438
$ws_LTW.Range("I132").Value = 3106.4119
$ws_LTW.Range("J132").Value = 4492.933
$ws_LTW.Range("K132").Value = 9319.235700000001
$ws_LTW.Range("L132").Value = 13478.799
$ws_LTW.Range("M132").Value = -6789.235700000001
$ws_LTW.Range("N132").Value = -18538.799

# WVR row 62
$ws_WVR.Range("H62").Value = 77975
$ws_WVR.Range("J62").Value = 77975
$ws_WVR.Range("L62").Value = 77975
$ws_WVR.Range("N62").Value = -79223

# WVR row 65
$ws_WVR.Range("H65").Value = 77975
$ws_WVR.Range("J65").Value = 77975
$ws_WVR.Range("L65").Value = 389875
$ws_WVR.Range("N65").Value = -396115

# WVR row 81
$ws_WVR.Range("H81").Value = 6538.5713
$ws_WVR.Range("I81").Value = 6811.6665
$ws_WVR.Range("J81").Value = 4900
$ws_WVR.Range("K81").Value = 13623.333
$ws_WVR.Range("L81").Value = 9800
$ws_WVR.Range("M81").Value = -12562.333
$ws_WVR.Range("N81").Value = -11922

# WVR row 84
$ws_WVR.Range("H84").Value = 6538.5713
$ws_WVR.Range("I84").Value = 6811.6665
$ws_WVR.Range("J84").Value = 4900
$ws_WVR.Range("K84").Value = 68116.66500000001
$ws_WVR.Range("L84").Value = 49000
$ws_WVR.Range("M84").Value = -62812.66500000001
$ws_WVR.Range("N84").Value = -59608

# WVR row 122
$ws_WVR.Range("H122").Value = 3545
$ws_WVR.Range("I122").Value = 1366.6666
$ws_WVR.Range("J122").Value = 3980.6667
$ws_WVR.Range("K122").Value = 4099.9998
$ws_WVR.Range("L122").Value = 11942.0001
$ws_WVR.Range("M122").Value = -1649.9998
$ws_WVR.Range("N122").Value = -16842.0001

# WVR row 132
$ws_WVR.Range("H132").Value = 5211162.5
$ws_WVR.Range("I132").Value = 2971.8096
$ws_WVR.Range("K132").Value = 8915.4288
$ws_WVR.Range("M132").Value = -6385.4288

